# Update the "F" column (想去人数 / want-to-go count) values on the
# "展览" and "全部类型" worksheets to reflect freshly generated stats.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 4272
$wsExhibit.Range("F3").Value = 2426
$wsExhibit.Range("F7").Value = 51
$wsExhibit.Range("F14").Value = 3294
$wsExhibit.Range("F15").Value = 223

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 4272
$wsAll.Range("F3").Value = 2426
$wsAll.Range("F8").Value = 51
$wsAll.Range("F18").Value = 3294
$wsAll.Range("F19").Value = 223
